$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = "[49.9354994553719, 50.165809478375394]"
$ws.Range("U2").Value = "[49.93791071998112, 50.1019518094902]"
$ws.Range("M3").Value = "[49.91134896963028, 50.16024335389203]"
$ws.Range("U3").Value = "[49.96063030804248, 50.129247661664536]"
$ws.Range("M4").Value = "[49.837995719056394, 50.1021289988351]"
$ws.Range("U4").Value = "[49.90057273592122, 50.06704102460903]"
$ws.Range("M5").Value = "[49.883285378044334, 50.16186632781965]"
$ws.Range("U5").Value = "[50.00757137491481, 50.173658580496]"
$ws.Range("M6").Value = "[49.93389387624619, 50.20839218938036]"
$ws.Range("U6").Value = "[49.95019896485398, 50.10819302958929]"
$ws.Range("M7").Value = "[49.80833588403445, 50.12761326237505]"
$ws.Range("U7").Value = "[49.89678835663548, 50.07411330504573]"
$ws.Range("M8").Value = "[49.80057059490395, 50.11970409578098]"
$ws.Range("U8").Value = "[49.89661567455387, 50.07157603886576]"
$ws.Range("M9").Value = "[49.9068960732158, 50.18341181546409]"
$ws.Range("U9").Value = "[49.8647632915018, 50.01211061462517]"
$ws.Range("M10").Value = "[49.859623601936555, 50.190194179870794]"
$ws.Range("U10").Value = "[49.86849986149844, 50.038412105444216]"
$ws.Range("M11").Value = "[49.917709881124054, 50.23120743290737]"
$ws.Range("U11").Value = "[49.93079442404706, 50.09425570322761]"
$ws.Range("M12").Value = "[49.99365937621855, 50.29618921106782]"
$ws.Range("U12").Value = "[49.93718527260484, 50.09615207362891]"
$ws.Range("M13").Value = "[49.94820190433965, 50.25817765616868]"
$ws.Range("U13").Value = "[49.883819072526045, 50.063932475508174]"
$ws.Range("M14").Value = "[49.93886659564557, 50.21161260269597]"
$ws.Range("U14").Value = "[49.932014766643704, 50.1209752382579]"
$ws.Range("M15").Value = "[49.735083035993824, 50.07578858282348]"
$ws.Range("U15").Value = "[49.85923124658646, 50.061608084228]"
$ws.Range("M16").Value = "[49.70022206898896, 50.073948021351235]"
$ws.Range("U16").Value = "[49.887886335433464, 50.08948466683381]"
